$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.984.45"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.540.11"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.16%  "
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0820"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.114"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "2.926.24"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "2.555.26"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.873"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "42.984.14"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.03%  "
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.67%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  +9.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0801"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.03%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.120"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "2.050.15"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").Value = "2.783.58"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.193"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.58%  "
